{"js": "// Replace the date line and every division-problem cell in the table with\n// the new values from the target edit. Each old string is unique in the\n// document, so a simple search+replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"2023-10-27 Friday\", \"2023-10-28 Saturday\"],\n  [\"47\u00f74=\", \"48\u00f77=\"],\n  [\"31\u00f77=\", \"72\u00f79=\"],\n  [\"77\u00f74=\", \"56\u00f73=\"],\n  [\"89\u00f75=\", \"39\u00f72=\"],\n  [\"13\u00f74=\", \"77\u00f72=\"],\n  [\"91\u00f78=\", \"55\u00f72=\"],\n  [\"35\u00f77=\", \"11\u00f79=\"],\n  [\"75\u00f72=\", \"74\u00f74=\"],\n  [\"72\u00f72=\", \"67\u00f72=\"],\n  [\"63\u00f72=\", \"14\u00f75=\"],\n  [\"71\u00f74=\", \"70\u00f77=\"],\n  [\"43\u00f77=\", \"26\u00f78=\"],\n  [\"56\u00f79=\", \"44\u00f72=\"],\n  [\"66\u00f77=\", \"25\u00f76=\"],\n  [\"86\u00f72=\", \"66\u00f79=\"],\n  [\"40\u00f72=\", \"71\u00f75=\"],\n  [\"40\u00f77=\", \"15\u00f75=\"],\n  [\"35\u00f72=\", \"29\u00f75=\"],\n  [\"43\u00f79=\", \"61\u00f74=\"],\n  [\"78\u00f73=\", \"57\u00f74=\"],\n  [\"42\u00f74=\", \"92\u00f76=\"],\n  [\"73\u00f77=\", \"62\u00f76=\"],\n  [\"69\u00f74=\", \"11\u00f73=\"],\n  [\"53\u00f77=\", \"67\u00f78=\"],\n  [\"60\u00f73=\", \"19\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every division-problem cell in the table with\n# the new values from the target edit. Each old string is unique in the\n# document, so a Find/Replace pass per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-10-27 Friday\", \"2023-10-28 Saturday\"),\n    @(\"47\u00f74=\", \"48\u00f77=\"),\n    @(\"31\u00f77=\", \"72\u00f79=\"),\n    @(\"77\u00f74=\", \"56\u00f73=\"),\n    @(\"89\u00f75=\", \"39\u00f72=\"),\n    @(\"13\u00f74=\", \"77\u00f72=\"),\n    @(\"91\u00f78=\", \"55\u00f72=\"),\n    @(\"35\u00f77=\", \"11\u00f79=\"),\n    @(\"75\u00f72=\", \"74\u00f74=\"),\n    @(\"72\u00f72=\", \"67\u00f72=\"),\n    @(\"63\u00f72=\", \"14\u00f75=\"),\n    @(\"71\u00f74=\", \"70\u00f77=\"),\n    @(\"43\u00f77=\", \"26\u00f78=\"),\n    @(\"56\u00f79=\", \"44\u00f72=\"),\n    @(\"66\u00f77=\", \"25\u00f76=\"),\n    @(\"86\u00f72=\", \"66\u00f79=\"),\n    @(\"40\u00f72=\", \"71\u00f75=\"),\n    @(\"40\u00f77=\", \"15\u00f75=\"),\n    @(\"35\u00f72=\", \"29\u00f75=\"),\n    @(\"43\u00f79=\", \"61\u00f74=\"),\n    @(\"78\u00f73=\", \"57\u00f74=\"),\n    @(\"42\u00f74=\", \"92\u00f76=\"),\n    @(\"73\u00f77=\", \"62\u00f76=\"),\n    @(\"69\u00f74=\", \"11\u00f73=\"),\n    @(\"53\u00f77=\", \"67\u00f78=\"),\n    @(\"60\u00f73=\", \"19\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
